$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1566.92
$ws.Range("I15").Value = 1566.92
$ws.Range("K15").Value = 4700.76
$ws.Range("M15").Value = -4531.76
$ws.Range("H18").Value = 2303
$ws.Range("I18").Value = 2303
$ws.Range("K18").Value = 2303
$ws.Range("M18").Value = -2019
$ws.Range("H32").Value = 6568
$ws.Range("I32").Value = 6292.875
$ws.Range("J32").Value = 7008.2
$ws.Range("K32").Value = 6292.875
$ws.Range("L32").Value = 7008.2
$ws.Range("M32").Value = -5966.875
$ws.Range("N32").Value = -7660.2
$ws.Range("H51").Value = 7649.5
$ws.Range("J51").Value = 8000
$ws.Range("L51").Value = 8000
$ws.Range("N51").Value = -8968
$ws.Range("H70").Value = 4573796
$ws.Range("J70").Value = 5561182
$ws.Range("L70").Value = 16683546
$ws.Range("N70").Value = -16684086
$ws.Range("H73").Value = 4573796
$ws.Range("J73").Value = 5561182
$ws.Range("L73").Value = 16683546
$ws.Range("N73").Value = -16685418
$ws.Range("H132").Value = 1620.7838
$ws.Range("I132").Value = 1559.4117
$ws.Range("J132").Value = 2316.3333
$ws.Range("K132").Value = 4678.2351
$ws.Range("L132").Value = 6948.999899999999
$ws.Range("M132").Value = -2148.2351
$ws.Range("N132").Value = -12008.9999
$ws.Range("H137").Value = 3308.4075
$ws.Range("I137").Value = 3179.1
$ws.Range("J137").Value = 3337.7954
$ws.Range("K137").Value = 9537.299999999999
$ws.Range("L137").Value = 10013.3862
$ws.Range("M137").Value = -6987.299999999999
$ws.Range("N137").Value = -15113.3862
$ws.Range("H138").Value = 2780.0435
$ws.Range("J138").Value = 3413.106
$ws.Range("L138").Value = 10239.318
$ws.Range("N138").Value = -20519.318
$ws.Range("H141").Value = 1993.1875
$ws.Range("I141").Value = 1782.6923
$ws.Range("K141").Value = 5348.0769
$ws.Range("M141").Value = -168.0769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2803.8306
$ws.Range("I32").Value = 2803.8306
$ws.Range("K32").Value = 2803.8306
$ws.Range("M32").Value = -2516.8306
$ws.Range("H61").Value = 6784.12
$ws.Range("I61").Value = 3170.9524
$ws.Range("J61").Value = 25753.25
$ws.Range("K61").Value = 3170.9524
$ws.Range("L61").Value = 25753.25
$ws.Range("M61").Value = -2958.9524
$ws.Range("N61").Value = -26177.25
$ws.Range("H102").Value = 3883.3333
$ws.Range("I102").Value = 3200
$ws.Range("K102").Value = 3200
$ws.Range("M102").Value = -1578
$ws.Range("H122").Value = 2183.5
$ws.Range("I122").Value = 1848.9412
$ws.Range("J122").Value = 2996
$ws.Range("K122").Value = 5546.8236
$ws.Range("L122").Value = 8988
$ws.Range("M122").Value = -3096.8236
$ws.Range("N122").Value = -13888
$ws.Range("H132").Value = 10012.963
$ws.Range("I132").Value = 8971.261
$ws.Range("K132").Value = 26913.783
$ws.Range("M132").Value = -24383.783
$ws.Range("H136").Value = 6784.12
$ws.Range("I136").Value = 3170.9524
$ws.Range("J136").Value = 25753.25
$ws.Range("K136").Value = 9512.8572
$ws.Range("L136").Value = 77259.75
$ws.Range("M136").Value = -6962.8572
$ws.Range("N136").Value = -82359.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 11259.434
$ws.Range("I105").Value = 11504.238
$ws.Range("K105").Value = 11504.238
$ws.Range("M105").Value = -9757.237999999999
$ws.Range("H107").Value = 1365.1666
$ws.Range("I107").Value = 1254.2941
$ws.Range("J107").Value = 3250
$ws.Range("K107").Value = 1254.2941
$ws.Range("L107").Value = 3250
$ws.Range("M107").Value = 665.7058999999999
$ws.Range("N107").Value = -7090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 914
$ws.Range("J2").Value = 250
$ws.Range("L2").Value = 250
$ws.Range("N2").Value = -476
$ws.Range("H31").Value = 29547.418
$ws.Range("I31").Value = 3631.36
$ws.Range("J31").Value = 65541.94500000001
$ws.Range("K31").Value = 3631.36
$ws.Range("L31").Value = 65541.94500000001
$ws.Range("M31").Value = -3336.36
$ws.Range("N31").Value = -66131.94500000001
$ws.Range("H34").Value = 29547.418
$ws.Range("I34").Value = 3631.36
$ws.Range("J34").Value = 65541.94500000001
$ws.Range("K34").Value = 3631.36
$ws.Range("L34").Value = 65541.94500000001
$ws.Range("M34").Value = -3429.36
$ws.Range("N34").Value = -65945.94500000001
$ws.Range("H105").Value = 6438.75
$ws.Range("I105").Value = 1750
$ws.Range("J105").Value = 8001.6665
$ws.Range("K105").Value = 1750
$ws.Range("L105").Value = 8001.6665
$ws.Range("M105").Value = -3
$ws.Range("N105").Value = -11495.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 125321.875
$ws.Range("J2").Value = 200254.8
$ws.Range("L2").Value = 1201528.8
$ws.Range("N2").Value = -1201754.8
$ws.Range("H12").Value = 43.142857
$ws.Range("J12").Value = 8.583333
$ws.Range("L12").Value = 25.749999
$ws.Range("N12").Value = -371.749999
$ws.Range("H33").Value = 400080.03
$ws.Range("I33").Value = 526392.7
$ws.Range("K33").Value = 3158356.2
$ws.Range("M33").Value = -3158073.2
$ws.Range("H50").Value = 66670576
$ws.Range("I50").Value = 83333470
$ws.Range("K50").Value = 250000410
$ws.Range("M50").Value = -249999929
$ws.Range("H53").Value = 66670576
$ws.Range("I53").Value = 83333470
$ws.Range("K53").Value = 250000410
$ws.Range("M53").Value = -249999929

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 8000567
$ws.Range("I11").Value = 11666666
$ws.Range("J11").Value = 4334468
$ws.Range("K11").Value = 11666666
$ws.Range("L11").Value = 4334468
$ws.Range("M11").Value = -11666527
$ws.Range("N11").Value = -4334746
$ws.Range("H24").Value = 53972.332
$ws.Range("I24").Value = 53972.332
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 53972.332
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -53799.332
$ws.Range("N24").ClearContents()
$ws.Range("H56").Value = 50000
$ws.Range("I56").Value = 50000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 50000
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("M56").Value = -49248
$ws.Range("H132").Value = 4661.3794
$ws.Range("I132").Value = 3703.3044
$ws.Range("J132").Value = 8334
$ws.Range("K132").Value = 11109.9132
$ws.Range("L132").Value = 25002
$ws.Range("M132").Value = -8579.913199999999
$ws.Range("N132").Value = -30062

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6151.846
$ws.Range("I7").Value = 4589.857
$ws.Range("K7").Value = 4589.857
$ws.Range("M7").Value = -4477.857
$ws.Range("H22").Value = 4056.9285
$ws.Range("I22").Value = 865.6667
$ws.Range("J22").Value = 6450.375
$ws.Range("K22").Value = 865.6667
$ws.Range("L22").Value = 6450.375
$ws.Range("M22").Value = -570.6667
$ws.Range("N22").Value = -7040.375
$ws.Range("H27").Value = 4056.9285
$ws.Range("I27").Value = 865.6667
$ws.Range("J27").Value = 6450.375
$ws.Range("K27").Value = 865.6667
$ws.Range("L27").Value = 6450.375
$ws.Range("M27").Value = -758.6667
$ws.Range("N27").Value = -6664.375
$ws.Range("H40").Value = 10297.55
$ws.Range("I40").Value = 11674.714
$ws.Range("K40").Value = 11674.714
$ws.Range("M40").Value = -11538.714
$ws.Range("H46").Value = 4115.3076
$ws.Range("I46").Value = 2916.1667
$ws.Range("J46").Value = 5143.143
$ws.Range("K46").Value = 2916.1667
$ws.Range("L46").Value = 5143.143
$ws.Range("M46").Value = -2728.1667
$ws.Range("N46").Value = -5519.143
$ws.Range("H55").Value = 1445.3636
$ws.Range("I55").Value = 431.7
$ws.Range("K55").Value = 431.7
$ws.Range("M55").Value = -258.7
$ws.Range("H68").Value = 3254.92
$ws.Range("I68").Value = 3320.8696
$ws.Range("K68").Value = 3320.8696
$ws.Range("M68").Value = -2571.8696
$ws.Range("H71").Value = 3254.92
$ws.Range("I71").Value = 3320.8696
$ws.Range("K71").Value = 16604.348
$ws.Range("M71").Value = -12860.348
$ws.Range("H100").Value = 5720.8
$ws.Range("I100").Value = 2532.6667
$ws.Range("J100").Value = 10503
$ws.Range("K100").Value = 2532.6667
$ws.Range("L100").Value = 10503
$ws.Range("M100").Value = -1991.6667
$ws.Range("N100").Value = -11585
$ws.Range("H126").Value = 6151.846
$ws.Range("I126").Value = 4589.857
$ws.Range("K126").Value = 13769.571
$ws.Range("M126").Value = -11299.571
$ws.Range("H132").Value = 6776.3
$ws.Range("I132").Value = 6251.737
$ws.Range("J132").Value = 7682.364
$ws.Range("K132").Value = 18755.211
$ws.Range("L132").Value = 23047.092
$ws.Range("M132").Value = -16225.211
$ws.Range("N132").Value = -28107.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 50000
$ws.Range("I26").Value = 50000
$ws.Range("K26").Value = 50000
$ws.Range("M26").Value = -49707
$ws.Range("H122").Value = 3060.8696
$ws.Range("I122").Value = 2373.7334
$ws.Range("J122").Value = 4349.25
$ws.Range("K122").Value = 7121.2002
$ws.Range("L122").Value = 13047.75
$ws.Range("M122").Value = -4671.2002
$ws.Range("N122").Value = -17947.75
